# Commit: "Added nudge for 10% expense reduction. Fixed previous bugs.
# PyFin seems to be working quite well, except that there seems to be a
# glitch with the pso."
#
# Semantic changes applied to Sheet1:
#   - Column B (Monthly IBT), rows 4-35: recomputed values (~10% "nudge"
#     reduction bug-fix cascade that also touches the later years).
#   - Column C (Monthly Expenses), rows 5-20: bug-fix values.
#   - Selection moved from H58 to K19.
#   - Workbook tab ratio nudged (cosmetic view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: Monthly IBT, rows 4-35 -------------------------------------
$colB = @{
    4  = 38750;  5  = 39750;  6  = 40750;  7  = 41750;  8  = 42750;
    9  = 43750;  10 = 44750;  11 = 45750;  12 = 46750;  13 = 47750;
    14 = 48750;  15 = 49750;  16 = 50750;  17 = 51750;  18 = 52750;
    19 = 53750;  20 = 54750;  21 = 55750;  22 = 56750;  23 = 57750;
    24 = 58750;  25 = 59750;  26 = 60750;  27 = 61750;  28 = 62750;
    29 = 63750;  30 = 64750;  31 = 65750;  32 = 66750;  33 = 67750;
    34 = 68750;  35 = 69750;
}
foreach ($row in $colB.Keys | Sort-Object) {
    $ws.Cells.Item($row, 2).Value = $colB[$row]
}

# --- Column C: Monthly Expenses, rows 5-20 --------------------------------
$colC = @{
    5  = 30001;  6  = 30002;  7  = 30003;  8  = 30004;  9  = 30005;
    10 = 30006;  11 = 30007;  12 = 30008;  13 = 30009;  14 = 30010;
    15 = 30011;  16 = 30012;  17 = 30013;  18 = 30014;  19 = 30015;
    20 = 30016;
}
foreach ($row in $colC.Keys | Sort-Object) {
    $ws.Cells.Item($row, 3).Value = $colC[$row]
}

# --- View state -------------------------------------------------------------
# Move the active selection from H58 to K19.
[void]$ws.Range("K19").Select()

# Tab ratio nudge (cosmetic; best-effort).
$excel.ActiveWindow.TabRatio = 0.5
